$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws3 = $wb.Worksheets.Item("table_v2")

# --- Update simulated-moment values (new calibration draw) on the "data" sheet ---
$ws1.Range("D5").Value = 2.373420435428225
$ws1.Range("E5").Value = 2.576630519628525
$ws1.Range("F5").Value = 0.08164715903196171
$ws1.Range("D6").Value = 0.07138971453727255
$ws1.Range("E6").Value = 0.06133794514834881
$ws1.Range("F6").Value = 0.002622678095011646
$ws1.Range("D7").Value = -0.3388462548494522
$ws1.Range("E7").Value = -0.3525642129778862
$ws1.Range("F7").Value = 0.01161199713255933
$ws1.Range("D8").Value = 0.2288157473547942
$ws1.Range("E8").Value = 0.2638800442814827
$ws1.Range("F8").Value = 0.008885795218818177
$ws1.Range("D9").Value = 2.493656150024008
$ws1.Range("E9").Value = 2.585287255048752
$ws1.Range("F9").Value = 0.08235110389677236
$ws1.Range("D10").Value = 0.2646393723224644
$ws1.Range("E10").Value = 0.2584711409658194
$ws1.Range("F10").Value = 0.01004802869067144
$ws1.Range("D11").Value = 2.423554827489142
$ws1.Range("E11").Value = 2.4094847646897
$ws1.Range("F11").Value = 0.07626731252655336
$ws1.Range("D12").Value = 0.309205812774586
$ws1.Range("E12").Value = 0.2851554799245871
$ws1.Range("F12").Value = 0.0121985665418822
$ws1.Range("D13").Value = 0.4515444406894221
$ws1.Range("E13").Value = 0.5315983342694378
$ws1.Range("F13").Value = 0.01916980495591796
$ws1.Range("D14").Value = 0.199905373436972
$ws1.Range("E14").Value = 0.1418508460156945
$ws1.Range("F14").Value = 0.007941641941221683
$ws1.Range("D15").Value = 0.1888404246655233
$ws1.Range("E15").Value = 0.2068904141768755
$ws1.Range("F15").Value = 0.01885783011450599
$ws1.Range("D16").Value = 0.2378275213560124
$ws1.Range("E16").Value = 0.2015302385505502
$ws1.Range("F16").Value = 0.01968563812959444
$ws1.Range("D17").Value = -0.1448787886555005
$ws1.Range("E17").Value = -0.05900479096240778
$ws1.Range("F17").Value = 0.01803628546388143
$ws1.Range("D18").Value = 0.00759650508292628
$ws1.Range("E18").Value = -0.02479339063389524
$ws1.Range("F18").Value = 0.01834683408094412
$ws1.Range("D19").Value = 0.4211344724013983
$ws1.Range("E19").Value = 0.3169023398360751
$ws1.Range("F19").Value = 0.01053797245042425
$ws1.Range("D20").Value = 0.6191579988240159
$ws1.Range("E20").Value = 0.5795648185578327
$ws1.Range("F20").Value = 0.02327032290503902
$ws1.Range("D21").Value = 0.1102882946930028
$ws1.Range("E21").Value = 0.3111629186435975
$ws1.Range("F21").Value = 0.01953587741858369
$ws1.Range("D22").Value = 0.02474953602968901
$ws1.Range("E22").Value = 0.03398225929140364
$ws1.Range("F22").Value = 0.006352793746508848

# --- Remove the two dropped moments ("Past portfolio/test and % expert"), rows 23-24 ---
# Deleting row 23 twice removes both rows and shifts row 25 (the SUM) up to row 23,
# auto-adjusting its formula from SUM(J5:J24) to SUM(J5:J22).
$ws1.Rows.Item(23).Delete()
$ws1.Rows.Item(23).Delete()

# --- View/selection bookkeeping ---
# Update the selection on table_v2 first (without leaving it as the active/selected tab),
# then make "data" the active sheet with its new selection - matches the target tabSelected
# + activeTab flip from table_v2 (index 2) back to data (index 0).
$ws3.Range("D3").Select()

$ws1.Activate()
$ws1.Range("C19").Select()
